$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 (existing "Conceptos clave: ..." slide, currently blank) gets its
# title + content filled in: "Conceptos clave: Versión"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$shpTitle3 = $s3.Shapes.Item(1)
$tr = $shpTitle3.TextFrame.TextRange
$tr.Text = "ón"
$tr.LanguageID = "es-419"
$rB = $tr.InsertBefore("Versi")
$rB.LanguageID = "en-US"
$rB = $tr.InsertBefore(": ")
$rB.LanguageID = "en-US"
$rB = $tr.InsertBefore("Conceptos clave")
$rB.LanguageID = "es-419"

$shpBody3 = $s3.Shapes.Item(2)
$trB = $shpBody3.TextFrame.TextRange
$trB.Text = " "
$trB.LanguageID = "es-419"
$lvl = $shpBody3.TextFrame.Ruler.Levels.Item(1)
$lvl.FirstMargin = 0
$lvl.LeftMargin = 0
$trB.ParagraphFormat.Bullet.Type = 0

# ---------------------------------------------------------------------------
# New slide: "Conceptos clave: Repositorio" (inserted right after slide 3)
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)

$tr = $s4.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Repositorio"
$tr.LanguageID = "en-US"
$rB = $tr.InsertBefore(": ")
$rB.LanguageID = "en-US"
$rB = $tr.InsertBefore("Conceptos clave")
$rB.LanguageID = "es-419"

$trB = $s4.Shapes.Item(2).TextFrame.TextRange
$trB.Text = " "
$trB.LanguageID = "es-419"
$lvl = $s4.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$lvl.FirstMargin = 0
$lvl.LeftMargin = 0
$trB.ParagraphFormat.Bullet.Type = 0

# ---------------------------------------------------------------------------
# New slide: "Conceptos clave: Control de Versiones"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)

$tr = $s5.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Versiones"
$tr.LanguageID = "en-US"
$rB = $tr.InsertBefore(": Control de ")
$rB.LanguageID = "en-US"
$rB = $tr.InsertBefore("Conceptos clave")
$rB.LanguageID = "es-419"

$trB = $s5.Shapes.Item(2).TextFrame.TextRange
$trB.Text = " "
$trB.LanguageID = "es-419"
$lvl = $s5.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$lvl.FirstMargin = 0
$lvl.LeftMargin = 0
$trB.ParagraphFormat.Bullet.Type = 0

# ---------------------------------------------------------------------------
# New slide: "Conceptos clave: Commit"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)

$tr = $s6.Shapes.Item(1).TextFrame.TextRange
$tr.Text = ": Commit"
$tr.LanguageID = "en-US"
$rB = $tr.InsertBefore("Conceptos clave")
$rB.LanguageID = "es-419"

$trB = $s6.Shapes.Item(2).TextFrame.TextRange
$trB.Text = " "
$trB.LanguageID = "es-419"
$lvl = $s6.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$lvl.FirstMargin = 0
$lvl.LeftMargin = 0
$trB.ParagraphFormat.Bullet.Type = 0

# ---------------------------------------------------------------------------
# New slide: "Conceptos clave: Push/Pull"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)

$tr = $s7.Shapes.Item(1).TextFrame.TextRange
$tr.Text = ": Push/Pull"
$tr.LanguageID = "en-US"
$rB = $tr.InsertBefore("Conceptos clave")
$rB.LanguageID = "es-419"

$trB = $s7.Shapes.Item(2).TextFrame.TextRange
$trB.Text = " "
$trB.LanguageID = "es-419"
$lvl = $s7.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$lvl.FirstMargin = 0
$lvl.LeftMargin = 0
$trB.ParagraphFormat.Bullet.Type = 0

# ---------------------------------------------------------------------------
# New slide: "Conceptos clave: Branch"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 2)

$tr = $s8.Shapes.Item(1).TextFrame.TextRange
$tr.Text = ": Branch"
$tr.LanguageID = "en-US"
$rB = $tr.InsertBefore("Conceptos clave")
$rB.LanguageID = "es-419"

$trB = $s8.Shapes.Item(2).TextFrame.TextRange
$trB.Text = " "
$trB.LanguageID = "es-419"
$lvl = $s8.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$lvl.FirstMargin = 0
$lvl.LeftMargin = 0
$trB.ParagraphFormat.Bullet.Type = 0

# ---------------------------------------------------------------------------
# New slide: "Conceptos clave: Merge/Conflict"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Add(9, 2)

$tr = $s9.Shapes.Item(1).TextFrame.TextRange
$tr.Text = ": Merge/Conflict"
$tr.LanguageID = "en-US"
$rB = $tr.InsertBefore("Conceptos clave")
$rB.LanguageID = "es-419"

$trB = $s9.Shapes.Item(2).TextFrame.TextRange
$trB.Text = " "
$trB.LanguageID = "es-419"
$lvl = $s9.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$lvl.FirstMargin = 0
$lvl.LeftMargin = 0
$trB.ParagraphFormat.Bullet.Type = 0
